$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E6"  = 16.139
    "C7"  = -12.917
    "B9"  = 5.907999999999999
    "C12" = -11.601
    "C14" = -12.845
    "E15" = 16.346
    "B18" = 5.01
    "B20" = 6.589999999999999
    "C26" = -13.279
    "B27" = 6.494
    "C27" = -13.059
    "C29" = -12.766
    "E33" = 17.437
    "B35" = 8.331000000000001
    "E35" = 16.527
    "C37" = -13.351
    "C38" = -12.986
    "E38" = 16.556
    "E43" = 17.038
    "E47" = 16.209
    "C51" = -11.471
    "E51" = 17.207
    "C52" = -11.621
    "C55" = -13.417
    "E57" = 16.394
    "E63" = 17.626
    "B69" = 5.425
    "C69" = -10.732
    "C70" = -11.959
    "E70" = 17.447
    "B76" = 6.723999999999999
    "B78" = 7.849999999999999
    "C81" = -13.926
    "B82" = 5.061
    "B83" = 5.111
    "C83" = -13.956
    "E88" = 16.235
    "B93" = 6.047000000000001
    "E99" = 16.488
    "C102" = -13.419
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
